$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.478.20"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").Value = "3.413.80"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'583.84"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").Value = "'180.20"
$ws.Range("E6").Value = "  +3.53%  "
$ws.Range("D7").Value = "'0.623"
$ws.Range("E7").Value = "  +5.14%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "3.412.64"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("E10").Value = "  +1.05%  "
$ws.Range("D11").Value = "'6.98"
$ws.Range("E11").Value = "  +1.97%  "
$ws.Range("D12").Value = "'0.413"
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("D13").Value = "4.014.29"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("D15").Value = "'29.14"
$ws.Range("E15").Value = "  -2.74%  "
$ws.Range("D16").Value = "66.496.61"
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("D18").Value = "3.414.17"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").Value = "'5.90"
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("D20").Value = "'13.82"
$ws.Range("E20").Value = "  +0.76%  "
$ws.Range("D21").Value = "'368.20"
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("D22").Value = "'7.59"
$ws.Range("E22").Value = "  -1.32%  "
$ws.Range("D23").Value = "'73.09"
$ws.Range("E23").Value = "  +1.87%  "
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").Value = "'0.0000125"
$ws.Range("E25").Value = "  +6.52%  "
$ws.Range("D26").Value = "'0.532"
$ws.Range("E26").Value = "  +1.06%  "
$ws.Range("D27").Value = "'9.81"
$ws.Range("E27").Value = "  +1.51%  "
$ws.Range("E28").Value = "  +1.38%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'1.99"
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'5.77"
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("D32").Value = "'23.30"
$ws.Range("E32").Value = "  -2.43%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").Value = "'7.05"
$ws.Range("D35").Value = "'1.26"
$ws.Range("E35").Value = "  -2.29%  "
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("D37").Value = "'163.19"
$ws.Range("E37").Value = "  +2.19%  "
$ws.Range("D38").Value = "'0.863"
$ws.Range("E38").Value = "  -1.39%  "
$ws.Range("D39").Value = "'27.51"
$ws.Range("E39").Value = "  -4.54%  "
$ws.Range("D40").Value = "'1.79"
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("D41").Value = "'2.63"
$ws.Range("E41").Value = "  +2.42%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.699.26"
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "'4.38"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("D45").Value = "'0.0686"
$ws.Range("E45").Value = "  +1.12%  "
$ws.Range("D46").Value = "'24.90"
$ws.Range("E46").Value = "  +3.42%  "
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").Value = "'39.89"
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").Value = "'333.71"
$ws.Range("E48").Value = "  +9.60%  "
$ws.Range("E49").Value = "  -0.88%  "
$ws.Range("E50").Value = "  +3.76%  "
$ws.Range("D51").Value = "'31.86"
$ws.Range("E51").Value = "  +5.37%  "